{"js": "// Update the date line and all twenty-five three-digit-by-one-digit\n// multiplication answers throughout the document's table cells.\n// Every \"old\" value below is unique in the document, so an exact,\n// case-sensitive whole-document search safely targets only the\n// intended run.\nconst replacements = [\n  [\"2025-05-24 Saturday\", \"2025-05-25 Sunday\"],\n  [\"744\u00d75=3720\", \"420\u00d77=2940\"],\n  [\"986\u00d77=6902\", \"965\u00d75=4825\"],\n  [\"876\u00d74=3504\", \"330\u00d74=1320\"],\n  [\"899\u00d74=3596\", \"459\u00d76=2754\"],\n  [\"570\u00d78=4560\", \"258\u00d73=774\"],\n  [\"382\u00d77=2674\", \"335\u00d75=1675\"],\n  [\"221\u00d74=884\", \"470\u00d77=3290\"],\n  [\"503\u00d75=2515\", \"462\u00d77=3234\"],\n  [\"889\u00d78=7112\", \"354\u00d79=3186\"],\n  [\"709\u00d76=4254\", \"336\u00d73=1008\"],\n  [\"269\u00d77=1883\", \"656\u00d79=5904\"],\n  [\"138\u00d73=414\", \"511\u00d74=2044\"],\n  [\"529\u00d73=1587\", \"495\u00d75=2475\"],\n  [\"663\u00d78=5304\", \"139\u00d72=278\"],\n  [\"277\u00d79=2493\", \"411\u00d79=3699\"],\n  [\"367\u00d72=734\", \"744\u00d73=2232\"],\n  [\"950\u00d75=4750\", \"888\u00d78=7104\"],\n  [\"451\u00d73=1353\", \"263\u00d78=2104\"],\n  [\"428\u00d73=1284\", \"447\u00d78=3576\"],\n  [\"399\u00d79=3591\", \"390\u00d77=2730\"],\n  [\"897\u00d73=2691\", \"177\u00d77=1239\"],\n  [\"386\u00d78=3088\", \"860\u00d74=3440\"],\n  [\"514\u00d74=2056\", \"909\u00d77=6363\"],\n  [\"198\u00d73=594\", \"656\u00d76=3936\"],\n  [\"209\u00d79=1881\", \"332\u00d79=2988\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all twenty-five three-digit-by-one-digit\n# multiplication answers throughout the document's table cells.\n# Every \"old\" value below is unique in the document, so a whole-document\n# Find/Replace (wdReplaceAll) for each exact string safely targets only\n# the intended run.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ old = \"2025-05-24 Saturday\"; new = \"2025-05-25 Sunday\" },\n    @{ old = \"744\u00d75=3720\"; new = \"420\u00d77=2940\" },\n    @{ old = \"986\u00d77=6902\"; new = \"965\u00d75=4825\" },\n    @{ old = \"876\u00d74=3504\"; new = \"330\u00d74=1320\" },\n    @{ old = \"899\u00d74=3596\"; new = \"459\u00d76=2754\" },\n    @{ old = \"570\u00d78=4560\"; new = \"258\u00d73=774\" },\n    @{ old = \"382\u00d77=2674\"; new = \"335\u00d75=1675\" },\n    @{ old = \"221\u00d74=884\"; new = \"470\u00d77=3290\" },\n    @{ old = \"503\u00d75=2515\"; new = \"462\u00d77=3234\" },\n    @{ old = \"889\u00d78=7112\"; new = \"354\u00d79=3186\" },\n    @{ old = \"709\u00d76=4254\"; new = \"336\u00d73=1008\" },\n    @{ old = \"269\u00d77=1883\"; new = \"656\u00d79=5904\" },\n    @{ old = \"138\u00d73=414\"; new = \"511\u00d74=2044\" },\n    @{ old = \"529\u00d73=1587\"; new = \"495\u00d75=2475\" },\n    @{ old = \"663\u00d78=5304\"; new = \"139\u00d72=278\" },\n    @{ old = \"277\u00d79=2493\"; new = \"411\u00d79=3699\" },\n    @{ old = \"367\u00d72=734\"; new = \"744\u00d73=2232\" },\n    @{ old = \"950\u00d75=4750\"; new = \"888\u00d78=7104\" },\n    @{ old = \"451\u00d73=1353\"; new = \"263\u00d78=2104\" },\n    @{ old = \"428\u00d73=1284\"; new = \"447\u00d78=3576\" },\n    @{ old = \"399\u00d79=3591\"; new = \"390\u00d77=2730\" },\n    @{ old = \"897\u00d73=2691\"; new = \"177\u00d77=1239\" },\n    @{ old = \"386\u00d78=3088\"; new = \"860\u00d74=3440\" },\n    @{ old = \"514\u00d74=2056\"; new = \"909\u00d77=6363\" },\n    @{ old = \"198\u00d73=594\"; new = \"656\u00d76=3936\" },\n    @{ old = \"209\u00d79=1881\"; new = \"332\u00d79=2988\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n}\n"}
